$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new header columns (Q, R) ---
$ws.Range("Q1").Value = "Người Ký"
$ws.Range("R1").Value = "Khách hàng"

# --- Fill the new columns for every data row (2-15) with the signer name ---
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 17).Value = "Trần Thị Hòa"
    $ws.Cells.Item($r, 18).Value = "Trần Thị Hòa"
}

# --- Re-fit the columns that now show wider/narrower content (mirrors the
#     "select all -> autofit column width" Excel performed after the edit) ---
$targetWidths = @(
    @("H", 12.022135416666666),
    @("I", 11.736979166666666),
    @("K", 10.451822916666666),
    @("L", 10.877604166666666),
    @("M", 21.022135416666668),
    @("N", 17.022135416666668),
    @("O", 21.022135416666668),
    @("P", 85.87760416666667),
    @("Q", 11.022135416666666)
)
foreach ($pair in $targetWidths) {
    $ws.Columns($pair[0]).ColumnWidth = $pair[1]
}

# --- Restore the selection Excel left behind after the edit ---
$ws.Range("P7").Select() | Out-Null
